$wb = $excel.ActiveWorkbook

# --- Update the descriptive text block on sheet "Hoja1" (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 9.45 = 38658.13 pesos
✅ 38658.13 pesos = 9.44 = 938.79 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 105.8
$wsTasas.Range("O10").Value = 4090.03
$wsTasas.Range("N12").Value = 4094
$wsTasas.Range("O12").Value = 99.42
